# Add columns I (I0) and J (IF) to the worksheet, mirroring the existing
# header style used by column H (bold, centered, bordered) and the plain
# numeric data pattern used by the other data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header cell (H1) onto the new
# header cells so they share the exact same style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Header text
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-10
$iValues = @(8, 9, 5, 7, 10, 7, 4, 3, 5)
$jValues = @(9, 9, 7, 8, 11, 9, 4, 3, 5)

for ($r = 0; $r -lt 9; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
